$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (02-11-2015, Wednesday) -> red fill, I6 = 1
$rng6 = $ws.Range("A6:J6")
$rng6.Interior.Color = 6184671
$ws.Cells.Item(6, 9).Value = 1

# Row 14 (02-19-2015) -> red fill, I14 = 1
$rng14 = $ws.Range("A14:J14")
$rng14.Interior.Color = 6184671
$ws.Cells.Item(14, 9).Value = 1

# Row 15 (02-20-2015, Friday) -> blue fill
$rng15 = $ws.Range("A15:J15")
$rng15.Interior.Color = 13411113

# B19 changes from shared string " " to boolean FALSE
$ws.Cells.Item(19, 2).Value = $false

# Fix formulas: remove redundant third FLOOR argument
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
